# The sheet used to hold a wide product-name/price table in A1:AD2.
# Replace it with a small two-column header: Name | Price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all of the old product/price data.
$ws.Range("A1:AD2").Clear()

# Write the new header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Price"
